{"js": "// The document originally contained the same 6-paragraph block of\n// \"placeholder\" Word content (5 text paragraphs + 1 trailing empty\n// paragraph) twice in a row. The edit removes the first copy of that\n// block, leaving only the second copy.\nconst body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = body.paragraphs.items;\nconst marker = \"Video provides a powerful way to help you prove your point.\";\n\n// Locate every paragraph that begins the duplicated block.\nconst matchIndexes = [];\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(marker) === 0) {\n    matchIndexes.push(i);\n  }\n}\n\n// Delete the paragraphs belonging to the first occurrence of the block.\n// That includes the first marker paragraph itself and every paragraph up\n// to (and including) the marker paragraph that starts the *second*\n// occurrence of the block (the duplicate block's own leading \"Video\n// provides...\" paragraph is removed too, since only one copy of the\n// whole block should remain). Fall back to deleting just the first\n// marker paragraph if there is no duplicate (defensive).\nconst startIndex = matchIndexes.length > 0 ? matchIndexes[0] : -1;\nconst endIndex = matchIndexes.length > 1 ? matchIndexes[1] : startIndex;\n\nif (startIndex !== -1) {\n  // Delete from the end backwards so earlier indexes stay valid.\n  for (let i = endIndex; i >= startIndex; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# The document originally contained the same 6-paragraph block of\n# \"placeholder\" Word content (5 text paragraphs + 1 trailing empty\n# paragraph) twice in a row. The edit removes the first copy of that\n# block, leaving only the second copy.\n$d = $word.ActiveDocument\n$marker = \"Video provides a powerful way to help you prove your point.\"\n\n$count = $d.Paragraphs.Count\n$matches = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t.StartsWith($marker)) {\n        $matches += $i\n    }\n}\n\nif ($matches.Count -gt 0) {\n    $startIndex = $matches[0]\n    # End just after the paragraph that starts the *second* occurrence of\n    # the block, so the duplicate block's own leading \"Video provides...\"\n    # paragraph is removed along with the whole first block.\n    if ($matches.Count -gt 1) {\n        $endIndex = $matches[1] + 1\n    } else {\n        $endIndex = $startIndex + 1\n    }\n\n    $startPos = $d.Paragraphs($startIndex).Range.Start\n    if ($endIndex -le $count) {\n        $endPos = $d.Paragraphs($endIndex).Range.Start\n    } else {\n        $endPos = $d.Content.End\n    }\n\n    $r = $d.Range($startPos, $endPos)\n    $r.Delete()\n}\n"}
